# Adds columns I (I0) and J (IF) with header labels and per-row values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Apply the same formatting used by the other header cells (e.g. H1) to I1/J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for I2:J40
$values = @(
    @(5, 5),
    @(3, 3),
    @(8, 8),
    @(6, 6),
    @(10, 10),
    @(1, 2),
    @(7, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(6, 7),
    @(8, 8),
    @(6, 7),
    @(6, 7),
    @(7, 7),
    @(7, 7),
    @(6, 7),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(5, 6),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(7, 8),
    @(6, 7),
    @(7, 8),
    @(4, 5),
    @(4, 5),
    @(7, 7),
    @(7, 7),
    @(5, 5),
    @(5, 6),
    @(3, 5),
    @(6, 6),
    @(5, 5),
    @(5, 5),
    @(6, 6)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
